# Regenerate the "K" (strikeouts) column (column G, formerly "Strike#")
# for each game row on Sheet1. The new values were recomputed upstream
# (regenerated std/mean + s_vals) and are written back into the sheet
# cell-by-cell, matching each row's game record (column A = game index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new K value (column G)
$kValues = @{
    2 = 0
    3 = 1
    4 = 2
    5 = 2
    6 = 1
    7 = 2
    8 = 1
    9 = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 2
    27 = 0
    28 = 2
    29 = 2
    30 = 2
    32 = 0
    33 = 3
    34 = 0
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 2
    40 = 0
    41 = 0
    42 = 1
    43 = 3
    44 = 0
    45 = 3
    46 = 0
    47 = 1
    48 = 0
    49 = 0
    50 = 0
    51 = 2
    52 = 1
    54 = 1
    56 = 2
    57 = 1
    58 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
